# Update automàtic: dades i banners [2026-02-28 04:20]
# Applies the per-cell text updates from the diff against
# src/data/resum_diari_meteocat.xlsx (sheet 'Dades_Meteo').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-28 04:18:33'
$ws.Range('E3').Value = '2026-02-28 04:18:35'
$ws.Range('H3').Value = "'83%"
$ws.Range('L3').Value = '20.9 km/h - 121º 3:57 TU'
$ws.Range('E4').Value = '2026-02-28 04:18:38'
$ws.Range('H4').Value = "'96%"
$ws.Range('O4').Value = '7.6 °C'
$ws.Range('E5').Value = '2026-02-28 04:18:40'
$ws.Range('N5').Value = '-0.4 °C 3:53 TU'
$ws.Range('O5').Value = '0.0 °C'
$ws.Range('E6').Value = '2026-02-28 04:18:43'
$ws.Range('O6').Value = '10.2 °C'
$ws.Range('E7').Value = '2026-02-28 04:18:46'
$ws.Range('H7').Value = "'76%"
$ws.Range('E8').Value = '2026-02-28 04:18:48'
$ws.Range('L8').Value = '16.6 km/h - 71º 3:48 TU'
$ws.Range('E9').Value = '2026-02-28 04:18:51'
$ws.Range('M9').Value = '8.5 °C 3:40 TU'
$ws.Range('O9').Value = '7.4 °C'
$ws.Range('E10').Value = '2026-02-28 04:18:53'
$ws.Range('O10').Value = '7.5 °C'
$ws.Range('E11').Value = '2026-02-28 04:18:56'
$ws.Range('H11').Value = "'92%"
$ws.Range('N11').Value = '2.9 °C 3:56 TU'
$ws.Range('O11').Value = '3.8 °C'
$ws.Range('E12').Value = '2026-02-28 04:18:58'
$ws.Range('O12').Value = '5.9 °C'
$ws.Range('E13').Value = '2026-02-28 04:19:01'
$ws.Range('H13').Value = "'85%"
$ws.Range('J13').Value = '1025.9 hPa'
$ws.Range('L13').Value = '11.5 km/h - 149º 3:48 TU'
$ws.Range('N13').Value = '-0.4 °C 3:32 TU'
$ws.Range('O13').Value = '1.7 °C'
$ws.Range('E14').Value = '2026-02-28 04:19:03'
$ws.Range('M14').Value = '10.8 °C 3:30 TU'
$ws.Range('O14').Value = '10.0 °C'
$ws.Range('E15').Value = '2026-02-28 04:19:06'
$ws.Range('O15').Value = '6.6 °C'
$ws.Range('E16').Value = '2026-02-28 04:19:08'
$ws.Range('H16').Value = "'63%"
$ws.Range('N16').Value = '-1.5 °C 3:49 TU'
$ws.Range('E17').Value = '2026-02-28 04:19:11'
$ws.Range('E18').Value = '2026-02-28 04:19:14'
$ws.Range('M18').Value = '9.0 °C 3:58 TU'
$ws.Range('O18').Value = '7.9 °C'
$ws.Range('E19').Value = '2026-02-28 04:19:16'
$ws.Range('E20').Value = '2026-02-28 04:19:19'
$ws.Range('H20').Value = "'36%"
$ws.Range('N20').Value = '-0.5 °C 3:41 TU'
$ws.Range('E21').Value = '2026-02-28 04:19:21'
$ws.Range('H21').Value = "'75%"
$ws.Range('N21').Value = '4.5 °C 3:59 TU'
$ws.Range('O21').Value = '5.6 °C'
$ws.Range('E22').Value = '2026-02-28 04:19:24'
$ws.Range('L22').Value = '22.3 km/h - 115º 3:59 TU'
$ws.Range('N22').Value = '-1.6 °C 3:30 TU'
$ws.Range('O22').Value = '-0.7 °C'
$ws.Range('E23').Value = '2026-02-28 04:19:27'
$ws.Range('H23').Value = "'69%"
$ws.Range('N23').Value = '-0.7 °C 3:35 TU'
$ws.Range('O23').Value = '0.1 °C'
$ws.Range('E24').Value = '2026-02-28 04:19:29'
$ws.Range('H24').Value = "'98%"
$ws.Range('J24').Value = '1023.3 hPa'
$ws.Range('N24').Value = '4.1 °C 3:48 TU'
$ws.Range('O24').Value = '6.2 °C'
$ws.Range('E25').Value = '2026-02-28 04:19:32'
$ws.Range('N25').Value = '-0.1 °C 3:59 TU'
$ws.Range('O25').Value = '1.2 °C'
$ws.Range('E26').Value = '2026-02-28 04:19:35'
$ws.Range('H26').Value = "'71%"
$ws.Range('J26').Value = '1024.0 hPa'
$ws.Range('O26').Value = '4.5 °C'
$ws.Range('E27').Value = '2026-02-28 04:19:37'
$ws.Range('N27').Value = '0.9 °C 3:51 TU'
$ws.Range('O27').Value = '2.7 °C'
$ws.Range('E28').Value = '2026-02-28 04:19:39'
$ws.Range('O28').Value = '6.5 °C'
$ws.Range('E29').Value = '2026-02-28 04:19:42'
$ws.Range('H29').Value = "'97%"
$ws.Range('L29').Value = '9.0 km/h - 16º 3:58 TU'
$ws.Range('M29').Value = '9.8 °C 3:53 TU'
$ws.Range('O29').Value = '8.5 °C'
$ws.Range('E30').Value = '2026-02-28 04:19:44'
$ws.Range('J30').Value = '1024.3 hPa'
$ws.Range('N30').Value = '7.3 °C 3:47 TU'
$ws.Range('E31').Value = '2026-02-28 04:19:47'
$ws.Range('H31').Value = "'94%"
$ws.Range('N31').Value = '9.7 °C 3:57 TU'
$ws.Range('O31').Value = '10.2 °C'
$ws.Range('E32').Value = '2026-02-28 04:19:49'
$ws.Range('H32').Value = "'87%"
$ws.Range('E33').Value = '2026-02-28 04:19:52'
$ws.Range('H33').Value = "'69%"
$ws.Range('L33').Value = '11.2 km/h - 67º 3:49 TU'
$ws.Range('E34').Value = '2026-02-28 04:19:54'
$ws.Range('H34').Value = "'72%"
$ws.Range('E35').Value = '2026-02-28 04:19:56'
$ws.Range('J35').Value = '1022.6 hPa'
$ws.Range('N35').Value = '6.3 °C 3:59 TU'
$ws.Range('O35').Value = '7.0 °C'
$ws.Range('E36').Value = '2026-02-28 04:19:59'
$ws.Range('L36').Value = '9.4 km/h - 275º 3:32 TU'
$ws.Range('E37').Value = '2026-02-28 04:20:02'
$ws.Range('N37').Value = '4.2 °C 3:53 TU'
$ws.Range('E38').Value = '2026-02-28 04:20:04'
$ws.Range('M38').Value = '9.3 °C 3:50 TU'
$ws.Range('O38').Value = '8.9 °C'
$ws.Range('E39').Value = '2026-02-28 04:20:06'
$ws.Range('N39').Value = '-1.0 °C 3:57 TU'
$ws.Range('O39').Value = '0.4 °C'
$ws.Range('E40').Value = '2026-02-28 04:20:09'
$ws.Range('N40').Value = '2.8 °C 3:59 TU'
$ws.Range('O40').Value = '3.7 °C'
$ws.Range('E41').Value = '2026-02-28 04:20:12'
$ws.Range('L41').Value = '12.2 km/h - 50º 3:56 TU'
$ws.Range('E42').Value = '2026-02-28 04:20:14'
$ws.Range('M42').Value = '8.3 °C 3:58 TU'
$ws.Range('E43').Value = '2026-02-28 04:20:17'
$ws.Range('H43').Value = "'83%"
$ws.Range('N43').Value = '3.5 °C 3:33 TU'
$ws.Range('O43').Value = '4.0 °C'
$ws.Range('E44').Value = '2026-02-28 04:20:19'
$ws.Range('H44').Value = "'91%"
$ws.Range('O44').Value = '-1.0 °C'
$ws.Range('E45').Value = '2026-02-28 04:20:22'
$ws.Range('H45').Value = "'89%"
$ws.Range('L45').Value = '11.9 km/h - 304º 3:36 TU'
$ws.Range('N45').Value = '6.4 °C 3:54 TU'
$ws.Range('O45').Value = '7.3 °C'
$ws.Range('E46').Value = '2026-02-28 04:20:24'
$ws.Range('J46').Value = '1023.0 hPa'
$ws.Range('N46').Value = '10.6 °C 3:39 TU'
